# Apply the data edit described by the commit "Read and write excel data by
# column name": refresh the header-ish first row of Sheet1 with new values,
# then update the active-cell selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Overwrite row 1 (columns A..J) with the new text values.
$ws.Cells.Item(1, 1).Value  = "DFSDF"
$ws.Cells.Item(1, 2).Value  = "SDFEWR"
$ws.Cells.Item(1, 3).Value  = "SDFWER"
$ws.Cells.Item(1, 4).Value  = "WERDSF"
$ws.Cells.Item(1, 5).Value  = "WEWW"
$ws.Cells.Item(1, 6).Value  = "DSFXCCXZ"
$ws.Cells.Item(1, 7).Value  = "ZXCSD"
$ws.Cells.Item(1, 8).Value  = "TRYTHG"
$ws.Cells.Item(1, 9).Value  = "DSFSDF"
$ws.Cells.Item(1, 10).Value = "DSFSDF"

# Move the active selection on Sheet1 from H17 to H13, matching the saved view.
$ws.Activate()
$ws.Range("H13").Select()
